$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32 (G32=5484)
$ws.Range("H32").Value = 2066.6667
$ws.Range("I32").Value = 2200
$ws.Range("J32").Value = 2000
$ws.Range("K32").Value = 2200
$ws.Range("L32").Value = 2000
$ws.Range("M32").Value = -1874
$ws.Range("N32").Value = -2652
# Row 138 (G138=44169)
$ws.Range("H138").Value = 2004495.1
$ws.Range("I138").Value = 2362.125
$ws.Range("J138").Value = 2946675.2
$ws.Range("K138").Value = 7086.375
$ws.Range("L138").Value = 8840025.600000001
$ws.Range("M138").Value = -1946.375
$ws.Range("N138").Value = -8850305.600000001
# Row 141 (G141=44161)
$ws.Range("H141").Value = 5164
$ws.Range("I141").Value = 4108.1113
$ws.Range("K141").Value = 12324.3339
$ws.Range("M141").Value = -7144.333899999998

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (G32=44147)
$ws.Range("H32").Value = 3576980.5
$ws.Range("I32").Value = 3707794.5
$ws.Range("K32").Value = 3707794.5
$ws.Range("M32").Value = -3707507.5
# Row 61 (G61=43999)
$ws.Range("H61").Value = 41669780
$ws.Range("I61").Value = 1363.375
$ws.Range("K61").Value = 1363.375
$ws.Range("M61").Value = -1151.375
# Row 106 (G106=18679)
$ws.Range("H106").Value = 38197.168
$ws.Range("J106").Value = 38197.168
$ws.Range("L106").Value = 38197.168
$ws.Range("N106").Value = -40721.168
# Row 123 (G123=34107)
$ws.Range("H123").Value = 73276.336
$ws.Range("J123").Value = 73276.336
$ws.Range("L123").Value = 73276.336
$ws.Range("N123").Value = -83076.336
# Row 136 (G136=43999)
$ws.Range("H136").Value = 41669780
$ws.Range("I136").Value = 1363.375
$ws.Range("K136").Value = 4090.125
$ws.Range("M136").Value = -1540.125

$ws = $wb.Worksheets.Item("BSM")
# Row 86 (G86=12526)
$ws.Range("H86").Value = 11410152
$ws.Range("J86").Value = 1967.2
$ws.Range("L86").Value = 1967.2
$ws.Range("N86").Value = -4213.2
# Row 89 (G89=12526)
$ws.Range("H89").Value = 11410152
$ws.Range("J89").Value = 1967.2
$ws.Range("L89").Value = 9836
$ws.Range("N89").Value = -21068
# Row 118 (G118=27137)
$ws.Range("H118").Value = 55000
$ws.Range("J118").Value = 55000
$ws.Range("L118").Value = 55000
$ws.Range("N118").Value = -58314

$ws = $wb.Worksheets.Item("CRP")
# Row 2 (G2=1820)
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").ClearContents()
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = 0
# Row 28 (G28=18348)
$ws.Range("H28").Value = 29388
$ws.Range("J28").Value = 29388
$ws.Range("L28").Value = 29388
$ws.Range("N28").Value = -29878
# Row 31 (G31=44023)
$ws.Range("H31").Value = 4950.3335
$ws.Range("I31").Value = 2531.4375
$ws.Range("J31").Value = 7226.9414
$ws.Range("K31").Value = 2531.4375
$ws.Range("L31").Value = 7226.9414
$ws.Range("M31").Value = -2236.4375
$ws.Range("N31").Value = -7816.9414
# Row 34 (G34=44023)
$ws.Range("H34").Value = 4950.3335
$ws.Range("I34").Value = 2531.4375
$ws.Range("J34").Value = 7226.9414
$ws.Range("K34").Value = 2531.4375
$ws.Range("L34").Value = 7226.9414
$ws.Range("M34").Value = -2329.4375
$ws.Range("N34").Value = -7630.9414
# Row 58 (G58=44021)
$ws.Range("H58").Value = 4796.241
$ws.Range("I58").Value = 2003.6364
$ws.Range("K58").Value = 2003.6364
$ws.Range("M58").Value = -1800.6364
# Row 86 (G86=12584)
$ws.Range("H86").Value = 44048064
$ws.Range("I86").Value = 37141900
$ws.Range("J86").Value = 55558340
$ws.Range("K86").Value = 37141900
$ws.Range("L86").Value = 55558340
$ws.Range("M86").Value = -37140777
$ws.Range("N86").Value = -55560586
# Row 89 (G89=12584)
$ws.Range("H89").Value = 44048064
$ws.Range("I89").Value = 37141900
$ws.Range("J89").Value = 55558340
$ws.Range("K89").Value = 185709500
$ws.Range("L89").Value = 277791700
$ws.Range("M89").Value = -185703884
$ws.Range("N89").Value = -277802932
# Row 105 (G105=19928)
$ws.Range("H105").Value = 7937755
$ws.Range("I105").Value = 10204901
$ws.Range("K105").Value = 10204901
$ws.Range("M105").Value = -10203154
# Row 125 (G125=34297)
$ws.Range("H125").Value = 47153
$ws.Range("J125").Value = 47153
$ws.Range("L125").Value = 47153
$ws.Range("N125").Value = -52073
# Row 132 (G132=44019)
$ws.Range("H132").Value = 3007.9185
$ws.Range("I132").Value = 2207.6128
$ws.Range("K132").Value = 6622.8384
$ws.Range("M132").Value = -4092.8384
# Row 136 (G136=44021)
$ws.Range("H136").Value = 4796.241
$ws.Range("I136").Value = 2003.6364
$ws.Range("K136").Value = 6010.9092
$ws.Range("M136").Value = -3460.9092
# Row 137 (G137=43231)
$ws.Range("H137").Value = 63000
$ws.Range("J137").Value = 69500
$ws.Range("L137").Value = 69500
$ws.Range("N137").Value = -79700

$ws = $wb.Worksheets.Item("CUL")
# Row 12 (G12=4854)
$ws.Range("H12").Value = 2941887.5
$ws.Range("I12").Value = 1866.8334
$ws.Range("J12").Value = 4545535
$ws.Range("K12").Value = 5600.5002
$ws.Range("L12").Value = 13636605
$ws.Range("M12").Value = -5427.5002
$ws.Range("N12").Value = -13636951
# Row 17 (G17=4640)
$ws.Range("H17").Value = 1114.8667
$ws.Range("I17").Value = 237.7
$ws.Range("J17").Value = 2869.2
$ws.Range("K17").Value = 713.0999999999999
$ws.Range("L17").Value = 8607.599999999999
$ws.Range("M17").Value = -544.0999999999999
$ws.Range("N17").Value = -8945.599999999999
# Row 23 (G23=4858)
$ws.Range("H23").Value = 466.125
$ws.Range("I23").Value = 339.66666
$ws.Range("J23").Value = 542
$ws.Range("K23").Value = 1018.99998
$ws.Range("L23").Value = 1626
$ws.Range("M23").Value = -783.9999799999999
$ws.Range("N23").Value = -2096
# Row 32 (G32=4731)
$ws.Range("H32").Value = 99.8
$ws.Range("I32").Value = 99
$ws.Range("K32").Value = 297
$ws.Range("M32").Value = -14
# Row 34 (G34=4749)
$ws.Range("H34").Value = 6217.3125
$ws.Range("J34").Value = 6599.1333
$ws.Range("L34").Value = 19797.3999
$ws.Range("N34").Value = -19965.3999
# Row 39 (G39=4712)
$ws.Range("H39").Value = 9444.546
$ws.Range("J39").Value = 9487.777
$ws.Range("L39").Value = 28463.331
$ws.Range("N39").Value = -29051.331
# Row 47 (G47=4663)
$ws.Range("H47").Value = 703
$ws.Range("I47").Value = 703
$ws.Range("K47").Value = 2109
$ws.Range("M47").Value = -1678
# Row 55 (G55=4733)
$ws.Range("H55").Value = 6259499.5
$ws.Range("J55").Value = 6259499.5
$ws.Range("L55").Value = 18778498.5
$ws.Range("N55").Value = -18778852.5
# Row 82 (G82=12856)
$ws.Range("H82").Value = 70000
$ws.Range("J82").Value = 70000
$ws.Range("L82").Value = 210000
$ws.Range("N82").Value = -210812
# Row 85 (G85=12856)
$ws.Range("H85").Value = 70000
$ws.Range("J85").Value = 70000
$ws.Range("L85").Value = 210000
$ws.Range("N85").Value = -212808
# Row 92 (G92=19841)
$ws.Range("H92").Value = 1328
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 1328
$ws.Range("K92").Value = 0
$ws.Range("L92").ClearContents()
$ws.Range("M92").Value = 3984
$ws.Range("N92").Value = -6480

$ws = $wb.Worksheets.Item("GSM")
# Row 113 (G113=27710)
$ws.Range("H113").Value = 6275.385
$ws.Range("I113").Value = 2620.1333
$ws.Range("J113").Value = 8559.916999999999
$ws.Range("K113").Value = 2620.1333
$ws.Range("L113").Value = 8559.916999999999
$ws.Range("M113").Value = -450.1333
$ws.Range("N113").Value = -12899.917
# Row 122 (G122=36182)
$ws.Range("H122").Value = 2133239.2
$ws.Range("I122").Value = 2899964.8
$ws.Range("J122").Value = 3445.889
$ws.Range("K122").Value = 8699894.399999999
$ws.Range("L122").Value = 10337.667
$ws.Range("M122").Value = -8697444.399999999
$ws.Range("N122").Value = -15237.667
# Row 132 (G132=44008)
$ws.Range("H132").Value = 3071.1875
$ws.Range("I132").Value = 1697.4445
$ws.Range("K132").Value = 5092.333500000001
$ws.Range("M132").Value = -2562.333500000001
# Row 136 (G136=42218)
$ws.Range("H136").Value = 24569.186
$ws.Range("J136").Value = 23624.695
$ws.Range("L136").Value = 70874.08499999999
$ws.Range("N136").Value = -75974.08499999999

$ws = $wb.Worksheets.Item("LTW")
# Row 2 (G2=2631)
$ws.Range("H2").Value = 32000
$ws.Range("J2").Value = 24000
$ws.Range("L2").Value = 24000
$ws.Range("N2").Value = -24224
# Row 7 (G7=36249)
$ws.Range("H7").Value = 5500.2
$ws.Range("I7").Value = 4466.6665
$ws.Range("J7").Value = 6345.8184
$ws.Range("K7").Value = 4466.6665
$ws.Range("L7").Value = 6345.8184
$ws.Range("M7").Value = -4354.6665
$ws.Range("N7").Value = -6569.8184
# Row 44 (G44=3658)
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").ClearContents()
$ws.Range("N44").Value = 0
# Row 46 (G46=5282)
$ws.Range("H46").Value = 2523.0312
$ws.Range("I46").Value = 2020.9474
$ws.Range("J46").Value = 3256.8462
$ws.Range("K46").Value = 2020.9474
$ws.Range("L46").Value = 3256.8462
$ws.Range("M46").Value = -1832.9474
$ws.Range("N46").Value = -3632.8462
# Row 57 (G57=4153)
$ws.Range("H57").Value = 5270.5
$ws.Range("I57").Value = 1541
$ws.Range("J57").Value = 9000
$ws.Range("K57").Value = 1541
$ws.Range("L57").Value = 9000
$ws.Range("M57").Value = -975
$ws.Range("N57").Value = -10132
# Row 126 (G126=36249)
$ws.Range("H126").Value = 5500.2
$ws.Range("I126").Value = 4466.6665
$ws.Range("J126").Value = 6345.8184
$ws.Range("K126").Value = 13399.9995
$ws.Range("L126").Value = 19037.4552
$ws.Range("M126").Value = -10929.9995
$ws.Range("N126").Value = -23977.4552
# Row 132 (G132=44058)
$ws.Range("H132").Value = 7819259
$ws.Range("I132").Value = 15627615
$ws.Range("J132").Value = 10903.719
$ws.Range("K132").Value = 46882845
$ws.Range("L132").Value = 32711.157
$ws.Range("M132").Value = -46880315
$ws.Range("N132").Value = -37771.157

$ws = $wb.Worksheets.Item("WVR")
# Row 81 (G81=12596)
$ws.Range("H81").Value = 30072556
$ws.Range("I81").Value = 1751315.5
$ws.Range("K81").Value = 3502631
$ws.Range("M81").Value = -3501570
# Row 84 (G84=12596)
$ws.Range("H84").Value = 30072556
$ws.Range("I84").Value = 1751315.5
$ws.Range("K84").Value = 17513155
$ws.Range("M84").Value = -17507851
# Row 132 (G132=44029)
$ws.Range("H132").Value = 5251.607
$ws.Range("I132").Value = 5128.591
$ws.Range("J132").Value = 5702.6665
$ws.Range("K132").Value = 15385.773
$ws.Range("L132").Value = 17107.9995
$ws.Range("M132").Value = -12855.773
$ws.Range("N132").Value = -22167.9995
# Row 133 (G133=41869)
$ws.Range("H133").Value = 148833.17
$ws.Range("J133").Value = 148833.17
$ws.Range("L133").Value = 148833.17
$ws.Range("N133").Value = -158953.17
